$d = $word.ActiveDocument

# 1. Update the owner e-mail domain from "example" to "gmail".
$d.Content.Find.Execute("example", $true, $false, $false, $false, $false, $true, 1, $false, "gmail", 2)

# 2. Fill in the kart track name inside the quotation marks.
$d.Content.Find.Execute("Ver as reservas recebidas para a pista “”.", $true, $false, $false, $false, $false, $true, 1, $false, "Ver as reservas recebidas para a pista “Kartódromo Internacional de Braga”.", 2)

# 3. Drop " para o próximo domingo" from the "Aceitar" bullet, stopping right before
#    the comma so its (non-bold) run keeps its own formatting untouched.
$full = $d.Content.Text
$needle = "Aceitar uma reserva para o próximo domingo"
$start = $full.IndexOf($needle)
$end = $start + $needle.Length
$aceitarRange = $d.Range($start, $end)
$aceitarRange.Text = "Aceitar uma reserva"

# 4. Underline the word "que" in the "Se tiver dúvidas" bullet.
$full = $d.Content.Text
$needle = "Se tiver dúvidas durante o processo, sinta-se à vontade para partilhar, mas não será dada ajuda a menos que peça explicitamente."
$start = $full.IndexOf($needle)
$queOffset = $needle.IndexOf(" que ")
$queStart = $start + $queOffset + 1
$queEnd = $queStart + 3
$queRange = $d.Range($queStart, $queEnd)
Write-Host ("que text: [" + $queRange.Text + "]")
$queRange.Font.Underline = 1

# 5. Mark the "Default Paragraph Font" style as semi-hidden (UI-gallery hidden-until-used).
$style = $d.Styles.Item("Tipodeletrapredefinidodopargrafo")
try {
    $style.Hidden = $true
} catch {
}
